$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "66.346.33"
$ws.Cells.Item(2, 5).Value = "  -0.70%  "

$ws.Cells.Item(3, 4).Value = "3.317.86"
$ws.Cells.Item(3, 5).Value = "  -1.78%  "

$ws.Cells.Item(4, 5).Value = "  -0.03%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "188.20"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +2.13%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "557.87"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.50%  "

$ws.Cells.Item(7, 5).Value = "  +0.08%  "

$ws.Cells.Item(8, 5).Value = "  -2.08%  "

$ws.Cells.Item(9, 4).Value = "3.311.66"
$ws.Cells.Item(9, 5).Value = "  -1.71%  "

$ws.Cells.Item(10, 5).Value = "  -1.75%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.587"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -1.83%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "47.64"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -1.19%  "

$ws.Cells.Item(13, 5).Value = "  +0.72%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "8.66"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -1.13%  "

$ws.Cells.Item(15, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(15, 4).Value = "3.847.42"
$ws.Cells.Item(15, 5).Value = "  -1.37%  "

$ws.Cells.Item(16, 2).Value = "BitcoinCash"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "629.28"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +3.05%  "

$ws.Cells.Item(17, 2).Value = "WrappedBTC"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(17, 4).Value = "66.352.78"
$ws.Cells.Item(17, 5).Value = "  -0.45%  "

$ws.Cells.Item(18, 2).Value = "Chainlink"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "18.10"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +1.68%  "

$ws.Cells.Item(19, 5).Value = "  -0.36%  "

$ws.Cells.Item(20, 4).Value = "3.318.89"
$ws.Cells.Item(20, 5).Value = "  -1.28%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "10.85"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -7.32%  "

$ws.Cells.Item(22, 5).Value = "  -0.82%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "18.16"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +6.39%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "103.06"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +6.21%  "

$ws.Cells.Item(25, 5).Value = "  -2.72%  "

$ws.Cells.Item(26, 5).Value = "  -3.41%  "

$ws.Cells.Item(27, 5).Value = "  +1.35%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.74"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -0.92%  "

$ws.Cells.Item(29, 5).Value = "  +0.63%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "8.71"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -2.04%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "30.41"
$ws.Cells.Item(31, 4).Style = "Normal"

$ws.Cells.Item(32, 5).Value = "  +3.66%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "6.39"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +0.06%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "560.12"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +2.51%  "

$ws.Cells.Item(35, 5).Value = "  -1.35%  "

$ws.Cells.Item(36, 4).Value = "3.864.95"
$ws.Cells.Item(36, 5).Value = "  +0.48%  "

$ws.Cells.Item(37, 5).Value = "  -0.19%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "58.16"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -1.33%  "

$ws.Cells.Item(39, 5).Value = "  +0.07%  "

$ws.Cells.Item(40, 4).Value = "0.0₃0734"
$ws.Cells.Item(40, 5).Value = "  +0.18%  "

$ws.Cells.Item(41, 2).Value = "Stacks"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "3.33"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -1.90%  "

$ws.Cells.Item(42, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "34.11"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +3.48%  "

$ws.Cells.Item(43, 2).Value = "Fetch.AI"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.72"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.50%  "

$ws.Cells.Item(44, 2).Value = "Kaspa"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.129"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.04%  "

$ws.Cells.Item(45, 5).Value = "  -4.88%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "3.21"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -15.84%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.0420"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -0.33%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "3.22"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +2.61%  "

$ws.Cells.Item(49, 5).Value = "  -1.21%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.59"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -3.54%  "

$ws.Cells.Item(51, 5).Value = "  +0.17%  "

